# Auto-generated script to update Titan_Profits market-data cells
# per scheduled-runner refresh (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 75726.21000000001
$ws.Range("I15").Value = 75726.21000000001
$ws.Range("K15").Value = 227178.63
$ws.Range("M15").Value = -227009.63
$ws.Range("H28").Value = 530271.1
$ws.Range("I28").Value = 794393.8
$ws.Range("J28").Value = 2025.7142
$ws.Range("K28").Value = 794393.8
$ws.Range("L28").Value = 2025.7142
$ws.Range("M28").Value = -793908.8
$ws.Range("N28").Value = -2995.7142
$ws.Range("H62").Value = 5570894.5
$ws.Range("I62").Value = 8555917
$ws.Range("J62").Value = 27281.428
$ws.Range("K62").Value = 8555917
$ws.Range("L62").Value = 27281.428
$ws.Range("M62").Value = -8555293
$ws.Range("N62").Value = -28529.428
$ws.Range("H65").Value = 5570894.5
$ws.Range("I65").Value = 8555917
$ws.Range("J65").Value = 27281.428
$ws.Range("K65").Value = 42779585
$ws.Range("L65").Value = 136407.14
$ws.Range("M65").Value = -42776465
$ws.Range("N65").Value = -142647.14
$ws.Range("H97").Value = 800
$ws.Range("I97").Value = 800
$ws.Range("K97").Value = 2400
$ws.Range("M97").Value = -1904
$ws.Range("H112").Value = 38962340
$ws.Range("J112").Value = 38962340
$ws.Range("L112").Value = 116887020
$ws.Range("N112").Value = -116889236
$ws.Range("H116").Value = 1771.3684
$ws.Range("I116").Value = 1953.3846
$ws.Range("K116").Value = 1953.3846
$ws.Range("M116").Value = 1488.6154
$ws.Range("H124").Value = 25390
$ws.Range("J124").Value = 32980
$ws.Range("L124").Value = 32980
$ws.Range("N124").Value = -42800
$ws.Range("H132").Value = 22312.674
$ws.Range("I132").Value = 28482.13
$ws.Range("K132").Value = 85446.39
$ws.Range("M132").Value = -82916.39
$ws.Range("H137").Value = 24392062
$ws.Range("I137").Value = 30304140
$ws.Range("J137").Value = 4743.75
$ws.Range("K137").Value = 90912420
$ws.Range("L137").Value = 14231.25
$ws.Range("M137").Value = -90909870
$ws.Range("N137").Value = -19331.25
$ws.Range("H138").Value = 5833870
$ws.Range("I138").Value = 2103329
$ws.Range("J138").Value = 7815720
$ws.Range("K138").Value = 6309987
$ws.Range("L138").Value = 23447160
$ws.Range("M138").Value = -6304847
$ws.Range("N138").Value = -23457440

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26774.361
$ws.Range("I32").Value = 7167.531
$ws.Range("J32").Value = 106835.586
$ws.Range("K32").Value = 7167.531
$ws.Range("L32").Value = 106835.586
$ws.Range("M32").Value = -6880.531
$ws.Range("N32").Value = -107409.586
$ws.Range("H45").Value = 1182.5
$ws.Range("I45").Value = 1184.5652
$ws.Range("J45").Value = 1166.6666
$ws.Range("K45").Value = 1184.5652
$ws.Range("L45").Value = 1166.6666
$ws.Range("M45").Value = -807.5652
$ws.Range("N45").Value = -1920.6666
$ws.Range("H74").Value = 4613.7026
$ws.Range("I74").Value = 917.72
$ws.Range("J74").Value = 12313.667
$ws.Range("K74").Value = 917.72
$ws.Range("L74").Value = 12313.667
$ws.Range("M74").Value = -43.72000000000003
$ws.Range("N74").Value = -14061.667
$ws.Range("H77").Value = 4613.7026
$ws.Range("I77").Value = 917.72
$ws.Range("J77").Value = 12313.667
$ws.Range("K77").Value = 4588.6
$ws.Range("L77").Value = 61568.335
$ws.Range("M77").Value = -220.6000000000004
$ws.Range("N77").Value = -70304.33499999999
$ws.Range("H132").Value = 3042.6
$ws.Range("I132").Value = 2281.4707
$ws.Range("J132").Value = 4660
$ws.Range("K132").Value = 6844.4121
$ws.Range("L132").Value = 13980
$ws.Range("M132").Value = -4314.4121
$ws.Range("N132").Value = -19040
$ws.Range("H139").Value = 48107.375
$ws.Range("J139").Value = 48107.375
$ws.Range("L139").Value = 48107.375
$ws.Range("N139").Value = -58387.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2220
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 2650
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 2650
$ws.Range("M7").Value = -387
$ws.Range("N7").Value = -2876
$ws.Range("H105").Value = 2766.0637
$ws.Range("I105").Value = 2611.9395
$ws.Range("J105").Value = 3129.3572
$ws.Range("K105").Value = 2611.9395
$ws.Range("L105").Value = 3129.3572
$ws.Range("M105").Value = -864.9395
$ws.Range("N105").Value = -6623.3572
$ws.Range("H134").Value = 2869.9714
$ws.Range("I134").Value = 2228.926
$ws.Range("J134").Value = 5033.5
$ws.Range("K134").Value = 6686.778
$ws.Range("L134").Value = 15100.5
$ws.Range("M134").Value = -4151.778
$ws.Range("N134").Value = -20170.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4223.43
$ws.Range("I31").Value = 2293.641
$ws.Range("J31").Value = 6104.975
$ws.Range("K31").Value = 2293.641
$ws.Range("L31").Value = 6104.975
$ws.Range("M31").Value = -1998.641
$ws.Range("N31").Value = -6694.975
$ws.Range("H34").Value = 4223.43
$ws.Range("I34").Value = 2293.641
$ws.Range("J34").Value = 6104.975
$ws.Range("K34").Value = 2293.641
$ws.Range("L34").Value = 6104.975
$ws.Range("M34").Value = -2091.641
$ws.Range("N34").Value = -6508.975
$ws.Range("H99").Value = 13891610
$ws.Range("I99").Value = 2474.9412
$ws.Range("J99").Value = 47622364
$ws.Range("K99").Value = 2474.9412
$ws.Range("L99").Value = 47622364
$ws.Range("M99").Value = -976.9412000000002
$ws.Range("N99").Value = -47625360
$ws.Range("H105").Value = 976.86206
$ws.Range("I105").Value = 829.3889
$ws.Range("J105").Value = 1218.1818
$ws.Range("K105").Value = 829.3889
$ws.Range("L105").Value = 1218.1818
$ws.Range("M105").Value = 917.6111
$ws.Range("N105").Value = -4712.1818
$ws.Range("H126").Value = 13891610
$ws.Range("I126").Value = 2474.9412
$ws.Range("J126").Value = 47622364
$ws.Range("K126").Value = 7424.823600000001
$ws.Range("L126").Value = 142867092
$ws.Range("M126").Value = -4954.823600000001
$ws.Range("N126").Value = -142872032
$ws.Range("H134").Value = 34093468
$ws.Range("I134").Value = 50001484
$ws.Range("J134").Value = 20836788
$ws.Range("K134").Value = 150004452
$ws.Range("L134").Value = 62510364
$ws.Range("M134").Value = -150001917
$ws.Range("N134").Value = -62515434

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 548
$ws.Range("I45").Value = 100
$ws.Range("J45").Value = 622.6667
$ws.Range("K45").Value = 300
$ws.Range("L45").Value = 1868.0001
$ws.Range("M45").Value = 232
$ws.Range("N45").Value = -2932.0001
$ws.Range("H98").Value = 164
$ws.Range("I98").Value = 234.33333
$ws.Range("J98").Value = 93.666664
$ws.Range("K98").Value = 702.99999
$ws.Range("L98").Value = 280.999992
$ws.Range("M98").Value = 795.00001
$ws.Range("N98").Value = -3276.999992
$ws.Range("H131").Value = 13335298
$ws.Range("J131").Value = 15153628
$ws.Range("L131").Value = 45460884
$ws.Range("N131").Value = -45470964
$ws.Range("H139").Value = 9561.200000000001
$ws.Range("I139").Value = 8906
$ws.Range("J139").Value = 9888.799999999999
$ws.Range("K139").Value = 26718
$ws.Range("L139").Value = 29666.4
$ws.Range("M139").Value = -21578
$ws.Range("N139").Value = -39946.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2350.9614
$ws.Range("I122").Value = 2351.1365
$ws.Range("J122").Value = 2350
$ws.Range("K122").Value = 7053.4095
$ws.Range("L122").Value = 7050
$ws.Range("M122").Value = -4603.4095
$ws.Range("N122").Value = -11950
$ws.Range("H126").Value = 3060
$ws.Range("I126").Value = 3171.4285
$ws.Range("K126").Value = 9514.2855
$ws.Range("M126").Value = -7044.2855
$ws.Range("H132").Value = 4722.3125
$ws.Range("I132").Value = 5665.3335
$ws.Range("J132").Value = 4156.5
$ws.Range("K132").Value = 16996.0005
$ws.Range("L132").Value = 12469.5
$ws.Range("M132").Value = -14466.0005
$ws.Range("N132").Value = -17529.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10136.1
$ws.Range("I136").Value = 9509.200000000001
$ws.Range("J136").Value = 10763
$ws.Range("K136").Value = 28527.6
$ws.Range("L136").Value = 32289
$ws.Range("M136").Value = -25977.6
$ws.Range("N136").Value = -37389

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3336834.2
$ws.Range("I81").Value = 10002001
$ws.Range("J81").Value = 4251
$ws.Range("K81").Value = 20004002
$ws.Range("L81").Value = 8502
$ws.Range("M81").Value = -20002941
$ws.Range("N81").Value = -10624
$ws.Range("H84").Value = 3336834.2
$ws.Range("I84").Value = 10002001
$ws.Range("J84").Value = 4251
$ws.Range("K84").Value = 100020010
$ws.Range("L84").Value = 42510
$ws.Range("M84").Value = -100014706
$ws.Range("N84").Value = -53118
$ws.Range("H122").Value = 1686.0555
$ws.Range("I122").Value = 1543.2667
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 4629.800099999999
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -2179.800099999999
$ws.Range("N122").Value = -12100
$ws.Range("H136").Value = 2515.3408
$ws.Range("I136").Value = 683.0294
$ws.Range("J136").Value = 8745.200000000001
$ws.Range("K136").Value = 2049.0882
$ws.Range("L136").Value = 26235.6
$ws.Range("M136").Value = 500.9117999999999
$ws.Range("N136").Value = -31335.6
